$p = $ppt.ActivePresentation

# --- Part 1: update the "1000 of 3-step ahead future paths" text boxes on every slide ---
$oldText = "1000 of 3-step ahead future paths"
$newText = "1000 3-steps ahead future paths"
$replaced = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    $stack = New-Object System.Collections.ArrayList
    [void]$stack.Add($s.Shapes)

    while ($stack.Count -gt 0) {
        $lastIdx = $stack.Count - 1
        $shapes = $stack[$lastIdx]
        $stack.RemoveAt($lastIdx)

        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.Type -eq 6) {
                [void]$stack.Add($shp.GroupItems)
            } elseif ($shp.HasTextFrame -eq -1) {
                $txt = $shp.TextFrame.TextRange.Text
                if ($txt -eq $oldText) {
                    $shp.TextFrame.TextRange.Text = $newText
                    $replaced = $replaced + 1
                }
            }
        }
    }
}
Write-Host ("Replaced text boxes: " + $replaced)

# --- Part 2: update the cached date field text ("16/05/2019" -> "20/05/2019")
#     on the slide master and every custom layout's Date placeholder ---
$oldDate = "16/05/2019"
$newDate = "20/05/2019"
$dateUpdates = 0

function Update-DatePlaceholder($shapes) {
    $count = 0
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -eq $true) {
            if ($shp.HasTextFrame -eq -1) {
                $cur = $shp.TextFrame.TextRange.Text
                if ($cur -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                    $count = $count + 1
                }
            }
        }
    }
    return $count
}

$m = $p.SlideMaster
$dateUpdates = $dateUpdates + (Update-DatePlaceholder $m.Shapes)

$layouts = $m.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $lyt = $layouts.Item($L)
    $dateUpdates = $dateUpdates + (Update-DatePlaceholder $lyt.Shapes)
}

Write-Host ("Date placeholders updated: " + $dateUpdates)
